# Daily attendance processing - 2025-11-19 21:18:48
# Swap the order of the two recorder names in the "Recorded By" (column G)
# cell for the rows below, reflecting the updated recording-order convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 4; Value = 'System, backup@backdoor.com' },
    @{ Row = 6; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 7; Value = 'admin@admin.com, System' },
    @{ Row = 10; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 11; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 12; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 13; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 14; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 15; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 17; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 18; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 19; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 20; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 21; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 22; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 24; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 26; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 29; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 30; Value = 'System, backup@backdoor.com' },
    @{ Row = 32; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 33; Value = 'admin@admin.com, System' },
    @{ Row = 36; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 37; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 38; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 39; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 40; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 41; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 43; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 44; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 45; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 46; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 47; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 48; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 50; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 52; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 55; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 56; Value = 'System, backup@backdoor.com' },
    @{ Row = 58; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 59; Value = 'admin@admin.com, System' },
    @{ Row = 62; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 63; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 64; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 65; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 66; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 67; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 69; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 70; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 71; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 72; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 73; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 74; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 76; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 78; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 83; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 84; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 85; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 86; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 87; Value = 'admin@admin.com, dnasr281@gmail.com' },
    @{ Row = 90; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 92; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 93; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 94; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 96; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 99; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 101; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 109; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 110; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 111; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 112; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 113; Value = 'admin@admin.com, dnasr281@gmail.com' },
    @{ Row = 116; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 118; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 119; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 120; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 122; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 125; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 127; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 135; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 136; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 137; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 138; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 139; Value = 'admin@admin.com, dnasr281@gmail.com' },
    @{ Row = 142; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 144; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 145; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 146; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 148; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 151; Value = 'dnasr281@gmail.com, System' },
    @{ Row = 153; Value = 'dnasr281@gmail.com, System' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
